$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 9
    4  = 3
    5  = 4
    6  = 5
    7  = 4
    8  = 2
    9  = 6
    10 = 5
    11 = 4
    12 = 4
    13 = 3
    14 = 4
    15 = 7
    16 = 9
    17 = 3
    18 = 5
    19 = 5
    20 = 2
    21 = 4
    22 = 6
    23 = 1
    24 = 5
    25 = 7
    26 = 8
    27 = 11
    28 = 13
    29 = 5
    30 = 5
    31 = 2
    32 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
